# Auto-generated edit script: updates market-data derived cells
# (currentAveragePrice / LevePrice / LeveProfit columns H-N) across all 8 job sheets,
# per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 2160.3333
$ws.Range("I20").Value = 740.5
$ws.Range("K20").Value = 740.5
$ws.Range("M20").Value = -510.5
$ws.Range("H33").Value = 590.7
$ws.Range("I33").Value = 154.8
$ws.Range("J33").Value = 1898.4
$ws.Range("K33").Value = 154.8
$ws.Range("L33").Value = 1898.4
$ws.Range("M33").Value = 74.19999999999999
$ws.Range("N33").Value = -2356.4
$ws.Range("H35").Value = 2160.3333
$ws.Range("I35").Value = 740.5
$ws.Range("K35").Value = 740.5
$ws.Range("M35").Value = -361.5
$ws.Range("H41").Value = 672.3333
$ws.Range("I41").Value = 370
$ws.Range("J41").Value = 732.8
$ws.Range("K41").Value = 370
$ws.Range("L41").Value = 732.8
$ws.Range("M41").Value = 70
$ws.Range("N41").Value = -1612.8

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2461.6428
$ws.Range("I110").Value = 1348.4
$ws.Range("K110").Value = 1348.4
$ws.Range("M110").Value = 696.5999999999999
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").Value = ""

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 306.2
$ws.Range("I11").Value = 167.5
$ws.Range("J11").Value = 398.66666
$ws.Range("K11").Value = 167.5
$ws.Range("L11").Value = 398.66666
$ws.Range("M11").Value = -27.5
$ws.Range("N11").Value = -678.66666
$ws.Range("H26").Value = 32985.5
$ws.Range("I26").Value = 32985.5
$ws.Range("K26").Value = 32985.5
$ws.Range("M26").Value = -32693.5
$ws.Range("H40").Value = 64500
$ws.Range("I40").Value = 64500
$ws.Range("K40").Value = 64500
$ws.Range("M40").Value = -64235

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 11020.1
$ws.Range("I4").Value = 10000.333
$ws.Range("K4").Value = 10000.333
$ws.Range("M4").Value = -9888.333000000001
$ws.Range("H15").Value = 4163.8
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 4163.8
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 4163.8
$ws.Range("M15").Value = ""
$ws.Range("N15").Value = -4503.8
$ws.Range("H94").Value = 3853.5454
$ws.Range("J94").Value = 3482.8
$ws.Range("L94").Value = 3482.8
$ws.Range("N94").Value = -4384.8

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 641.3
$ws.Range("I98").Value = 632.5
$ws.Range("K98").Value = 1897.5
$ws.Range("M98").Value = -399.5
$ws.Range("H109").Value = 955.75
$ws.Range("I109").Value = 807.6667
$ws.Range("J109").Value = 1400
$ws.Range("K109").Value = 2423.0001
$ws.Range("L109").Value = 4200
$ws.Range("M109").Value = -1383.0001
$ws.Range("N109").Value = -6280
$ws.Range("H112").Value = 44610.11
$ws.Range("I112").Value = 1499
$ws.Range("J112").Value = 49999
$ws.Range("K112").Value = 4497
$ws.Range("L112").Value = 149997
$ws.Range("M112").Value = -3389
$ws.Range("N112").Value = -152213

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 23929.334
$ws.Range("J15").Value = 23929.334
$ws.Range("L15").Value = 23929.334
$ws.Range("N15").Value = -24505.334
$ws.Range("H80").Value = 2902.889
$ws.Range("I80").Value = 2995.75
$ws.Range("J80").Value = 2828.6
$ws.Range("K80").Value = 2995.75
$ws.Range("L80").Value = 2828.6
$ws.Range("M80").Value = -1997.75
$ws.Range("N80").Value = -4824.6
$ws.Range("H81").Value = 23929.334
$ws.Range("J81").Value = 23929.334
$ws.Range("L81").Value = 23929.334
$ws.Range("N81").Value = -25925.334
$ws.Range("H83").Value = 2902.889
$ws.Range("I83").Value = 2995.75
$ws.Range("J83").Value = 2828.6
$ws.Range("K83").Value = 14978.75
$ws.Range("L83").Value = 14143
$ws.Range("M83").Value = -9986.75
$ws.Range("N83").Value = -24127
$ws.Range("H84").Value = 23929.334
$ws.Range("J84").Value = 23929.334
$ws.Range("L84").Value = 71788.00199999999
$ws.Range("N84").Value = -81772.00199999999
$ws.Range("H102").Value = 957.5
$ws.Range("I102").Value = 951.4286
$ws.Range("K102").Value = 951.4286
$ws.Range("M102").Value = 670.5714
$ws.Range("H113").Value = 2177.7778
$ws.Range("I113").Value = 1514.2858
$ws.Range("J113").Value = 4500
$ws.Range("K113").Value = 1514.2858
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = 655.7141999999999
$ws.Range("N113").Value = -8840
$ws.Range("H122").Value = 4640.5713
$ws.Range("I122").Value = 4817.909
$ws.Range("J122").Value = 3990.3333
$ws.Range("K122").Value = 14453.727
$ws.Range("L122").Value = 11970.9999
$ws.Range("M122").Value = -12003.727
$ws.Range("N122").Value = -16870.9999
$ws.Range("H126").Value = 250001220
$ws.Range("I126").Value = 250001220
$ws.Range("K126").Value = 750003660
$ws.Range("M126").Value = -750001190
$ws.Range("H132").Value = 3299
$ws.Range("I132").Value = 2840.5715
$ws.Range("K132").Value = 8521.7145
$ws.Range("M132").Value = -5991.7145

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2648.3333
$ws.Range("I7").Value = 2473
$ws.Range("J7").Value = 2999
$ws.Range("K7").Value = 2473
$ws.Range("L7").Value = 2999
$ws.Range("M7").Value = -2361
$ws.Range("N7").Value = -3223
$ws.Range("H22").Value = 3539
$ws.Range("I22").Value = 2350.5
$ws.Range("J22").Value = 3878.5715
$ws.Range("K22").Value = 2350.5
$ws.Range("L22").Value = 3878.5715
$ws.Range("M22").Value = -2055.5
$ws.Range("N22").Value = -4468.5715
$ws.Range("H27").Value = 3539
$ws.Range("I27").Value = 2350.5
$ws.Range("J27").Value = 3878.5715
$ws.Range("K27").Value = 2350.5
$ws.Range("L27").Value = 3878.5715
$ws.Range("M27").Value = -2243.5
$ws.Range("N27").Value = -4092.5715
$ws.Range("H61").Value = 2413
$ws.Range("I61").Value = 1821
$ws.Range("K61").Value = 1821
$ws.Range("M61").Value = -1619
$ws.Range("H80").Value = 24684.5
$ws.Range("J80").Value = 24684.5
$ws.Range("L80").Value = 24684.5
$ws.Range("N80").Value = -26930.5
$ws.Range("H83").Value = 24684.5
$ws.Range("J83").Value = 24684.5
$ws.Range("L83").Value = 74053.5
$ws.Range("N83").Value = -85285.5
$ws.Range("H92").Value = 22111
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = ""
$ws.Range("H96").Value = 45000
$ws.Range("J96").Value = 45000
$ws.Range("L96").Value = 45000
$ws.Range("N96").Value = -50492
$ws.Range("H113").Value = 2413
$ws.Range("I113").Value = 1821
$ws.Range("K113").Value = 1821
$ws.Range("M113").Value = 349
$ws.Range("H126").Value = 2648.3333
$ws.Range("I126").Value = 2473
$ws.Range("J126").Value = 2999
$ws.Range("K126").Value = 7419
$ws.Range("L126").Value = 8997
$ws.Range("M126").Value = -4949
$ws.Range("N126").Value = -13937

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1815.1538
$ws.Range("I126").Value = 1556.8572
$ws.Range("J126").Value = 2116.5
$ws.Range("K126").Value = 4670.571599999999
$ws.Range("L126").Value = 6349.5
$ws.Range("M126").Value = -2200.571599999999
$ws.Range("N126").Value = -11289.5
$ws.Range("H136").Value = 2262.6086
$ws.Range("I136").Value = 1811.75
$ws.Range("K136").Value = 5435.25
$ws.Range("M136").Value = -2885.25
$ws.Range("H137").Value = 105000
$ws.Range("J137").Value = 105000
$ws.Range("L137").Value = 105000
$ws.Range("N137").Value = -115200

